# Rework cursed dremora items (#504)
# Insert a new "Dremora" material row at row 14 of the Weapons sheet,
# pushing the existing rows 14-38 down to 15-39.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weapons")

# Insert a new row before the current row 14 - this shifts rows 14..38
# down to 15..39 and keeps all of their existing data/formatting intact.
$ws.Rows.Item(14).EntireRow.Insert()

# Populate the newly inserted row 14 with the new Dremora entry.
$ws.Range("A14").Value = "Dremora"
$ws.Range("B14").Value = 7
$ws.Range("C14").Value = 7
$ws.Range("D14").Value = 10

# Match the final selection left behind by the edit.
[void]$ws.Range("H14").Select()
